$d = $word.ActiveDocument

# Find the paragraph containing "LOB1004: Cálculo II (Requisito fraco)" and
# remove the three paragraphs that follow it:
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "© 2020 ..." footer paragraph
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*LOB1004: Cálculo II (Requisito fraco)*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # Delete the three paragraphs immediately following the target paragraph.
    $d.Paragraphs.Item($target + 1).Range.Delete()
    $d.Paragraphs.Item($target + 1).Range.Delete()
    $d.Paragraphs.Item($target + 1).Range.Delete()
}

Write-Host "done"
